# Fix duplicate name issue: swap the data for rows 8 and 9 on the
# "Female_25m" sheet so that "Sara Alonso Lopez" (Stavanger, 31.10.2021,
# 29,97) is listed before "Sanna Josefin Husan Ehrnholm" (Trondheim,
# 18.06.2016, 29,96).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Female_25m")

# Columns A, B, D, E contain the values that differ between the two rows;
# C (Poeng), F (Pool) and G (Gender) are identical for both rows, but we
# swap the full row (A:G) for robustness.
$row8 = @(
    $ws.Cells.Item(8, 1).Value2,
    $ws.Cells.Item(8, 2).Value2,
    $ws.Cells.Item(8, 3).Value2,
    $ws.Cells.Item(8, 4).Value2,
    $ws.Cells.Item(8, 5).Value2,
    $ws.Cells.Item(8, 6).Value2,
    $ws.Cells.Item(8, 7).Value2
)

$row9 = @(
    $ws.Cells.Item(9, 1).Value2,
    $ws.Cells.Item(9, 2).Value2,
    $ws.Cells.Item(9, 3).Value2,
    $ws.Cells.Item(9, 4).Value2,
    $ws.Cells.Item(9, 5).Value2,
    $ws.Cells.Item(9, 6).Value2,
    $ws.Cells.Item(9, 7).Value2
)

for ($col = 1; $col -le 7; $col++) {
    $ws.Cells.Item(8, $col).Value = $row9[$col - 1]
    $ws.Cells.Item(9, $col).Value = $row8[$col - 1]
}
